$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$oldVol = "1"
$newVol = "2"
$volPos = $volText.LastIndexOf($oldVol)
$volCell.Characters($volPos + 1, $oldVol.Length).Text = $newVol

$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value2
$oldStart = "12/30/2024"
$newStart = "1/6/2025"
$startPos = $weekText.IndexOf($oldStart)
$weekCell.Characters($startPos + 1, $oldStart.Length).Text = $newStart

$weekText2 = $weekCell.Value2
$oldEnd = "1/5/2025"
$newEnd = "1/12/2025"
$endPos = $weekText2.IndexOf($oldEnd)
$weekCell.Characters($endPos + 1, $oldEnd.Length).Text = $newEnd

# --- Data table updates ---

$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("F15").Value = 3
$ws.Range("C16").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 2
$ws.Range("K17").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 233.333333333333
$ws.Range("I16").Value = 5
$ws.Range("C16").Copy($ws.Range("J16"))
$ws.Range("J16").Value = 2
$ws.Range("K17").Copy($ws.Range("K16"))
$ws.Range("K16").Value = 150
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -16.666666666666
$ws.Range("N16").Value = -84.848484848484
$ws.Range("C14").Copy($ws.Range("C17"))
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = -75
$ws.Range("J17").Value = 4
$ws.Range("K17").Copy($ws.Range("L17"))
$ws.Range("L17").Value = -100
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = -28.571428571428
$ws.Range("N18").Value = -90
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 12
$ws.Range("J19").Value = 13
$ws.Range("K19").Value = -7.692307692307
$ws.Range("M19").Value = -25
$ws.Range("N19").Value = -52
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 500
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = -8.333333333333
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 350
$ws.Range("L20").Value = 80
$ws.Range("M20").Value = 125
$ws.Range("N20").Value = -94.078947368421
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 18.75
$ws.Range("F21").Value = 72
$ws.Range("H21").Value = 1.408450704225
$ws.Range("I21").Value = 31
$ws.Range("J21").Value = 27
$ws.Range("K21").Value = 14.814814814814
$ws.Range("L21").Value = 24
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -88.076923076923
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C16").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K17").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("C16").Copy($ws.Range("J22"))
$ws.Range("J22").Value = 1
$ws.Range("K17").Copy($ws.Range("K22"))
$ws.Range("K22").Value = 0
$ws.Range("K17").Copy($ws.Range("L22"))
$ws.Range("L22").Value = 0
$ws.Range("K17").Copy($ws.Range("M22"))
$ws.Range("M22").Value = -50
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 48
$ws.Range("F24").Value = 159
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 40.707964601769
$ws.Range("I24").Value = 60
$ws.Range("J24").Value = 38
$ws.Range("K24").Value = 57.894736842105
$ws.Range("L24").Value = 39.534883720930
$ws.Range("M24").Value = 160.869565217391
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = 38.095238095238
$ws.Range("F25").Value = 116
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = 34.883720930232
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 31
$ws.Range("K25").Value = 54.838709677419
$ws.Range("L25").Value = 65.517241379310
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 16
$ws.Range("H26").Value = -33.333333333333
$ws.Range("C16").Copy($ws.Range("I26"))
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 20
$ws.Range("L26").Value = -33.333333333333
$ws.Range("M26").Value = -14.285714285714
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("F27").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("J28").Value = 2
$ws.Range("K17").Copy($ws.Range("L28"))
$ws.Range("L28").Value = -100
$ws.Range("J44").Value = 454
$ws.Range("K44").Value = -32.035928143712
$ws.Range("L44").Value = -38.731443994601
$ws.Range("M44").Value = -55.836575875486
$ws.Range("N44").Value = -66.119402985074
$ws.Range("J46").Value = 943
$ws.Range("K46").Value = -51.740020470829
$ws.Range("L46").Value = -64.226100151745
$ws.Range("M46").Value = -85.714285714285
$ws.Range("N46").Value = -87.166575939031